$wb = $excel.ActiveWorkbook

# --- Numeric value updates (re-run / simulation refresh) ---
$ws = $wb.Worksheets.Item("Combined_Data")
$ws.Range("I80").Value = 0.11315
$ws.Range("I83").Value = 0.341728
$ws.Range("I85").Value = 386.934
$ws.Range("I86").Value = 0.248265
$ws.Range("G87").Value = 0.77532
$ws.Range("H87").Value = 87.2929
$ws.Range("I87").Value = 146.397
$ws.Range("J87").Value = 5460
$ws.Range("G88").Value = 1.16284
$ws.Range("H88").Value = 96.1178
$ws.Range("I88").Value = 837.028
$ws.Range("J88").Value = 8189
$ws.Range("I89").Value = 0.936518
$ws.Range("G90").Value = 0.7288559999999999
$ws.Range("H90").Value = 86.9903
$ws.Range("I90").Value = 1.28894
$ws.Range("J90").Value = 954
$ws.Range("G91").Value = 2.75269
$ws.Range("H91").Value = 89.99169999999999
$ws.Range("I91").Value = 574.172
$ws.Range("J91").Value = 3603
$ws.Range("I92").Value = 0.485469
$ws.Range("G93").Value = 0.399162
$ws.Range("H93").Value = 93.4579
$ws.Range("I93").Value = 543.3869999999999
$ws.Range("J93").Value = 2811
$ws.Range("G94").Value = 0.440058
$ws.Range("H94").Value = 98.5308
$ws.Range("I94").Value = 748.5549999999999
$ws.Range("J94").Value = 3099
$ws.Range("I95").Value = 1.34131
$ws.Range("G96").Value = 0.696004
$ws.Range("H96").Value = 87.5767
$ws.Range("I96").Value = 112.863
$ws.Range("J96").Value = 911
$ws.Range("G97").Value = 0.738788
$ws.Range("H97").Value = 97.3139
$ws.Range("I97").Value = 340.181
$ws.Range("J97").Value = 967
$ws.Range("I104").Value = 0.35478
$ws.Range("G105").Value = 5.52849
$ws.Range("H105").Value = 9.39071
$ws.Range("I105").Value = 711.768
$ws.Range("J105").Value = 38933
$ws.Range("G106").Value = 5.76548
$ws.Range("H106").Value = 80.7516
$ws.Range("I106").Value = 3153.3
$ws.Range("J106").Value = 40602
$ws.Range("I107").Value = 1.14079
$ws.Range("G108").Value = 5.59401
$ws.Range("H108").Value = 0.150007
$ws.Range("I108").Value = 234.004
$ws.Range("J108").Value = 7322
$ws.Range("G109").Value = 10.7472
$ws.Range("H109").Value = 60.925
$ws.Range("I109").Value = 4128.21
$ws.Range("J109").Value = 14067

$ws = $wb.Worksheets.Item("Speed_Throughput")
$ws.Range("C4").Value = 4.7705
$ws.Range("D4").Value = 7.0378
$ws.Range("C5").Value = 0.2214
$ws.Range("D5").Value = 0.5264

$ws = $wb.Worksheets.Item("Speed_Packet_Loss")
$ws.Range("C4").Value = 43.9313
$ws.Range("D4").Value = 45.9823
$ws.Range("C5").Value = 96.81
$ws.Range("D5").Value = 5.251

$ws = $wb.Worksheets.Item("Speed_Avg_Delay")
$ws.Range("C4").Value = 450.028
$ws.Range("D4").Value = 957.7024
$ws.Range("C5").Value = 91.8579
$ws.Range("D5").Value = 222.4556

$ws = $wb.Worksheets.Item("Interferers_Throughput")
$ws.Range("C4").Value = 2.7275
$ws.Range("D4").Value = 5.8309
$ws.Range("C5").Value = 2.2644
$ws.Range("D5").Value = 5.1323

$ws = $wb.Worksheets.Item("Interferers_Packet_Loss")
$ws.Range("C4").Value = 69.3882
$ws.Range("D4").Value = 42.611
$ws.Range("C5").Value = 71.3532
$ws.Range("D5").Value = 41.9796

$ws = $wb.Worksheets.Item("Interferers_Avg_Delay")
$ws.Range("C4").Value = 204.2019
$ws.Range("D4").Value = 487.1066
$ws.Range("C5").Value = 337.684
$ws.Range("D5").Value = 886.6904

$ws = $wb.Worksheets.Item("Packet_Size_Throughput")
$ws.Range("C4").Value = 1.8325
$ws.Range("D4").Value = 3.3174
$ws.Range("C5").Value = 3.1594
$ws.Range("D5").Value = 6.9671

$ws = $wb.Worksheets.Item("Packet_Size_Packet_Loss")
$ws.Range("C4").Value = 73.24850000000001
$ws.Range("D4").Value = 40.2238
$ws.Range("C5").Value = 67.4928
$ws.Range("D5").Value = 44.1019

$ws = $wb.Worksheets.Item("Packet_Size_Avg_Delay")
$ws.Range("C4").Value = 334.3009
$ws.Range("D4").Value = 722.5688
$ws.Range("C5").Value = 207.585
$ws.Range("D5").Value = 708.7259

$ws = $wb.Worksheets.Item("Traffic_Rate_Throughput")
$ws.Range("C5").Value = 2.0293
$ws.Range("D5").Value = 2.7112
$ws.Range("C7").Value = 5.2698
$ws.Range("D7").Value = 8.4123

$ws = $wb.Worksheets.Item("Traffic_Rate_Packet_Loss")
$ws.Range("C5").Value = 65.2929
$ws.Range("D5").Value = 46.3284
$ws.Range("C7").Value = 81.31950000000001
$ws.Range("D7").Value = 30.3591

$ws = $wb.Worksheets.Item("Traffic_Rate_Avg_Delay")
$ws.Range("C5").Value = 88.0381
$ws.Range("D5").Value = 180.6893
$ws.Range("C6").Value = 0.2588
$ws.Range("D6").Value = 0.3861
$ws.Range("C7").Value = 724.532
$ws.Range("D7").Value = 1101.241

$ws = $wb.Worksheets.Item("SmartV3_Data")
$ws.Range("I8").Value = 0.11315
$ws.Range("I11").Value = 0.341728
$ws.Range("I13").Value = 386.934
$ws.Range("I14").Value = 0.248265
$ws.Range("G15").Value = 0.77532
$ws.Range("H15").Value = 87.2929
$ws.Range("I15").Value = 146.397
$ws.Range("J15").Value = 5460
$ws.Range("G16").Value = 1.16284
$ws.Range("H16").Value = 96.1178
$ws.Range("I16").Value = 837.028
$ws.Range("J16").Value = 8189
$ws.Range("I17").Value = 0.936518
$ws.Range("G18").Value = 0.7288559999999999
$ws.Range("H18").Value = 86.9903
$ws.Range("I18").Value = 1.28894
$ws.Range("J18").Value = 954
$ws.Range("G19").Value = 2.75269
$ws.Range("H19").Value = 89.99169999999999
$ws.Range("I19").Value = 574.172
$ws.Range("J19").Value = 3603
$ws.Range("I20").Value = 0.485469
$ws.Range("G21").Value = 0.399162
$ws.Range("H21").Value = 93.4579
$ws.Range("I21").Value = 543.3869999999999
$ws.Range("J21").Value = 2811
$ws.Range("G22").Value = 0.440058
$ws.Range("H22").Value = 98.5308
$ws.Range("I22").Value = 748.5549999999999
$ws.Range("J22").Value = 3099
$ws.Range("I23").Value = 1.34131
$ws.Range("G24").Value = 0.696004
$ws.Range("H24").Value = 87.5767
$ws.Range("I24").Value = 112.863
$ws.Range("J24").Value = 911
$ws.Range("G25").Value = 0.738788
$ws.Range("H25").Value = 97.3139
$ws.Range("I25").Value = 340.181
$ws.Range("J25").Value = 967
$ws.Range("I32").Value = 0.35478
$ws.Range("G33").Value = 5.52849
$ws.Range("H33").Value = 9.39071
$ws.Range("I33").Value = 711.768
$ws.Range("J33").Value = 38933
$ws.Range("G34").Value = 5.76548
$ws.Range("H34").Value = 80.7516
$ws.Range("I34").Value = 3153.3
$ws.Range("J34").Value = 40602
$ws.Range("I35").Value = 1.14079
$ws.Range("G36").Value = 5.59401
$ws.Range("H36").Value = 0.150007
$ws.Range("I36").Value = 234.004
$ws.Range("J36").Value = 7322
$ws.Range("G37").Value = 10.7472
$ws.Range("H37").Value = 60.925
$ws.Range("I37").Value = 4128.21
$ws.Range("J37").Value = 14067

$ws = $wb.Worksheets.Item("Overall_Throughput")
$ws.Range("B3").Value = 2.4959
$ws.Range("C3").Value = 5.4589

$ws = $wb.Worksheets.Item("Overall_Packet_Loss")
$ws.Range("B3").Value = 70.3707
$ws.Range("C3").Value = 42.0092
$ws.Range("D3").Value = 99.2654

$ws = $wb.Worksheets.Item("Overall_Avg_Delay")
$ws.Range("B3").Value = 270.9429
$ws.Range("C3").Value = 713.4814
$ws.Range("D3").Value = 0.0566
$ws.Range("F3").Value = 4128.21

$ws = $wb.Worksheets.Item("Distance_Throughput")
$ws.Range("C5").Value = 4.3631
$ws.Range("D5").Value = 7.1414
$ws.Range("C6").Value = 3.1247
$ws.Range("D6").Value = 5.5172

$ws = $wb.Worksheets.Item("Distance_Packet_Loss")
$ws.Range("C5").Value = 51.3948
$ws.Range("D5").Value = 43.4938
$ws.Range("C6").Value = 59.7172
$ws.Range("D6").Value = 46.3231

$ws = $wb.Worksheets.Item("Distance_Avg_Delay")
$ws.Range("C5").Value = 349.1941
$ws.Range("D5").Value = 534.4529
$ws.Range("C6").Value = 463.6347
$ws.Range("D6").Value = 1078.9149

# --- Column width auto-fit side effects (content width changed) ---
$ws = $wb.Worksheets.Item("Overall_Packet_Loss")
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666

$ws = $wb.Worksheets.Item("Speed_Avg_Delay")
$ws.Columns.Item(4).ColumnWidth = 9.166666666666666

$ws = $wb.Worksheets.Item("Traffic_Rate_Avg_Delay")
$ws.Columns.Item(4).ColumnWidth = 9.166666666666666

# --- Analysis metadata timestamp update ---
$ws = $wb.Worksheets.Item("Analysis_Summary")
$ws.Range("B6").Value = "2025-08-17 12:32:34"

